$d = $word.ActiveDocument

# The paragraph ends with a hyperlink followed immediately by the
# "_GoBack" bookmark (bookmarkStart/bookmarkEnd), then the paragraph
# mark. We want the four new runs inserted between the hyperlink and
# the bookmark, i.e. the bookmark should end up AFTER the new text
# (it stays "anchored" right before the paragraph's very end once the
# new text has been typed there).
#
# Appending via InsertAfter at the paragraph end reliably produces
# clean runs with no inherited run formatting, but it always lands
# the new text after the existing bookmark. So: first append the
# clean text (landing after the bookmark), then relocate the
# "_GoBack" bookmark back to its original spot -- immediately after
# the hyperlink and before the newly inserted text.

function Insert-AtParaEnd($doc, $text) {
    $p = $doc.Paragraphs($doc.Paragraphs.Count)
    $insertAt = $p.Range.End - 1
    $r = $doc.Range($insertAt, $insertAt)
    $r.InsertAfter($text)
}

# Locate the insertion point: right after the hyperlink's visible
# text, before the bookmark, in the *original* document layout.
$hyperlink = $d.Hyperlinks(1)
$bookmarkPos = $hyperlink.Range.End

Insert-AtParaEnd $d ".  Both"
Insert-AtParaEnd $d " the"
Insert-AtParaEnd $d " command line and the Team Explorer GUI built into Visual Studio were used for this"
Insert-AtParaEnd $d " in order to get more experience with both. "

# Re-anchor "_GoBack" immediately after the hyperlink (before the new
# runs we just appended), matching its original relative position.
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
